$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "44.647.09"
Set-TextValue "E2" "  +3.82%  "
Set-TextValue "D3" "2.433.13"
Set-TextValue "E3" "  +2.43%  "
Set-TextValue "E4" "  -0.07%  "
Set-TextValue "D5" "311.97"
Set-TextValue "E5" "  +3.31%  "
Set-TextValue "D6" "102.13"
Set-TextValue "E6" "  +6.12%  "
Set-TextValue "E7" "  +1.97%  "
Set-TextValue "E8" "  -0.04%  "
Set-TextValue "E9" "  +2.54%  "
Set-TextValue "D10" "35.53"
Set-TextValue "E10" "  +3.94%  "
Set-TextValue "E11" "  +1.69%  "
Set-TextValue "E12" "  +1.19%  "
Set-TextValue "D13" "18.81"
Set-TextValue "E13" "  +3.30%  "
Set-TextValue "E14" "  +3.04%  "
Set-TextValue "D15" "2.813.43"
Set-TextValue "E15" "  +2.46%  "
Set-TextValue "D16" "2.454.42"
Set-TextValue "E16" "  +2.68%  "
Set-TextValue "E17" "  +4.49%  "
Set-TextValue "D18" "44.565.12"
Set-TextValue "E18" "  +3.73%  "
Set-TextValue "D19" "12.56"
Set-TextValue "E19" "  +2.98%  "
Set-TextValue "E20" "  +1.82%  "
Set-TextValue "D21" "0.0₃0910"
Set-TextValue "E21" "  +2.50%  "
Set-TextValue "D22" "68.96"
Set-TextValue "E22" "  +1.26%  "
Set-TextValue "D23" "2.32"
Set-TextValue "E23" "  +4.04%  "
Set-TextValue "D24" "241.33"
Set-TextValue "E24" "  +2.70%  "
Set-TextValue "E25" "  +1.61%  "
Set-TextValue "E26" "  +0.01%  "
Set-TextValue "D27" "25.32"
Set-TextValue "E27" "  +1.81%  "
Set-TextValue "E28" "  -4.20%  "
Set-TextValue "D29" "9.69"
Set-TextValue "E29" "  +4.68%  "
Set-TextValue "D30" "33.54"
Set-TextValue "E30" "  +6.37%  "
Set-TextValue "E31" "  +16.52%  "
Set-TextValue "D32" "19.53"
Set-TextValue "E32" "  +11.42%  "
Set-TextValue "E33" "  +2.82%  "
Set-TextValue "D35" "0.0761"
Set-TextValue "E35" "  +3.91%  "
Set-TextValue "E36" "  +2.99%  "
Set-TextValue "D37" "4.53"
Set-TextValue "E37" "  +3.98%  "
Set-TextValue "E38" "  +4.17%  "
Set-TextValue "D39" "126.67"
Set-TextValue "E39" "  +8.44%  "
Set-TextValue "B40" "Stellar"
Set-TextValue "C40" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D40" "0.109"
Set-TextValue "E40" "  +0.90%  "
Set-TextValue "B41" "WEMIXToken"
Set-TextValue "C41" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D41" "2.20"
Set-TextValue "E41" "  -4.62%  "
Set-TextValue "D42" "21.97"
Set-TextValue "E42" "  -0.14%  "
Set-TextValue "E43" "  +3.73%  "
Set-TextValue "D44" "1.950.27"
Set-TextValue "E44" "  +0.31%  "
Set-TextValue "D45" "2.17"
Set-TextValue "E45" "  +2.20%  "
Set-TextValue "E46" "  +8.68%  "
Set-TextValue "D47" "9.76"
Set-TextValue "D48" "1.69"
Set-TextValue "E48" "  +11.34%  "
Set-TextValue "B49" "MultiversX"
Set-TextValue "C49" "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextValue "D49" "53.54"
Set-TextValue "E49" "  +2.91%  "
Set-TextValue "B50" "BitcoinSV"
Set-TextValue "C50" "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextValue "D50" "74.02"
Set-TextValue "E50" "  +2.75%  "
Set-TextValue "B51" "TrustWalletToken"
Set-TextValue "C51" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D51" "1.16"
Set-TextValue "E51" "  +1.81%  "
